# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the repulled dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    3  = 1
    4  = -1
    5  = 2
    6  = -1
    8  = 3
    9  = 1
    11 = 3
    12 = 0
    13 = 1
    14 = -5
    16 = -5
    17 = 6
    18 = 1
    19 = 1
    20 = -2
    22 = -2
    24 = 4
    25 = 1
    26 = 2
    27 = 1
    28 = 6
    29 = -1
    30 = -3
    31 = -3
    32 = -5
    33 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
